# "updated shops form for excel"
# Adds three new header columns (pincode, village, tehsil) after the
# existing operatorName/address headers, matching the style already used
# by the header row (A1/B1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header cell's formatting (12pt font, vertical-center)
# onto the new header cells before/while filling in their text.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("C1:E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("C1").Value = "pincode"
$ws.Range("D1").Value = "village"
$ws.Range("E1").Value = "tehsil"

$ws.Range("F3").Select() | Out-Null
